# Update cryptos list values (price / volume) per latest scrape.
# Rows 19-22 and 42-43 also had their Coin/Link values swapped (row reordering
# reflected as value changes at fixed row positions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.273.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.885.28'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.33'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4671'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2834'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06586'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.75'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07775'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '98.11'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.888.31'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.119'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6687'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '283.89'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +10.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.269.41'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.135.32'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.60'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.363'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.74%  '
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007307'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.44%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.172'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.339'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.16'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.995'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.05%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09738'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.465'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.482'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.177'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04697'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7093'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.095'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01871'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.660'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.527'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.37'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.32%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8700'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.09%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.972'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '104.03'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.13%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4204'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '987.90'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.228'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.276'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1162'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.00%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.13%  '
